# Fruta / hortaliza, semanal
# Adds this week's price block (date 2022-03-17, serial 44637) for
# "Agrícola del Norte S.A. de Arica - Piña" at the top of the data table
# (row 129), pushing the previously existing rows (129-164) down by one
# block of 4 rows (to 133-168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows (one full "quality" block: Especial/Primera/Segunda/Tercera)
# right above the current first data block for this market/product, shifting
# all the existing rows below it down by 4.
$ws.Range("A129:T132").Insert()

# Common (constant) values shared by every row in this block.
$mercadoId = 1
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$fecha = 44637
$codreg = 15
$tipo = "Fruta"
$productoId = 100108
$producto = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria = "Piña"
$variedad = "Caramelo"
$origen = "Ecuador"

# Per-quality data for the new block: Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg, Kg/unidad
$rows = @(
    @{ Row = 129; Calidad = "Especial"; Volumen = 200; PrecioMin = 18000; PrecioMax = 19000; PrecioProm = 18500; PrecioKg = 1850; KgUnidad = 10 },
    @{ Row = 130; Calidad = "Primera";  Volumen = 250; PrecioMin = 18000; PrecioMax = 19000; PrecioProm = 18500; PrecioKg = 1542; KgUnidad = 12 },
    @{ Row = 131; Calidad = "Segunda";  Volumen = 270; PrecioMin = 18000; PrecioMax = 19000; PrecioProm = 18500; PrecioKg = 1321; KgUnidad = 14 },
    @{ Row = 132; Calidad = "Tercera";  Volumen = 200; PrecioMin = 18000; PrecioMax = 19000; PrecioProm = 18500; PrecioKg = 1156; KgUnidad = 16 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PrecioMin
    $ws.Cells.Item($row, 15).Value = $r.PrecioMax
    $ws.Cells.Item($row, 16).Value = $r.PrecioProm
    $ws.Cells.Item($row, 17).Value = "`$/caja $($r.KgUnidad) unidades"
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
